$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3286
$ws.Range("E2").Value = 68
$ws.Range("F2").Value = 68
$ws.Range("G2").Value = -133
$ws.Range("H2").Value = -114
$ws.Range("I2").Value = -116
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 2875
$ws.Range("L2").Value = 1886
$ws.Range("M2").Value = 989
$ws.Range("N2").Value = 909
$ws.Range("O2").Value = 80
$ws.Range("P2").Value = 140
$ws.Range("Q2").Value = 125
$ws.Range("R2").Value = -115
$ws.Range("S2").Value = -15
$ws.Range("T2").Value = 109
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 1734
$ws.Range("W2").Value = 2.06
$ws.Range("X2").Value = -3.47
$ws.Range("Y2").Value = -11.97
$ws.Range("Z2").Value = -3.9
$ws.Range("AA2").Value = 190.67
$ws.Range("AB2").Value = 554.35
$ws.Range("AC2").Value = -353
$ws.Range("AD2").Value = -2.82
$ws.Range("AE2").Value = 2846
$ws.Range("AF2").Value = 0.35
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 32960505

# Row 3
$ws.Range("D3").Value = 2562
$ws.Range("E3").Value = -110
$ws.Range("F3").Value = -110
$ws.Range("G3").Value = -659
$ws.Range("H3").Value = -567
$ws.Range("I3").Value = -532
$ws.Range("J3").Value = -34
$ws.Range("K3").Value = 2450
$ws.Range("L3").Value = 1783
$ws.Range("M3").Value = 667
$ws.Range("N3").Value = 620
$ws.Range("O3").Value = 47
$ws.Range("P3").Value = 140
$ws.Range("Q3").Value = 124
$ws.Range("R3").Value = 12
$ws.Range("S3").Value = -66
$ws.Range("T3").Value = 14
$ws.Range("U3").Value = 110
$ws.Range("V3").Value = 1644
$ws.Range("W3").Value = -4.29
$ws.Range("X3").Value = -22.11
$ws.Range("Y3").Value = -69.62
$ws.Range("Z3").Value = -21.28
$ws.Range("AA3").Value = 267.33
$ws.Range("AB3").Value = 173.69
$ws.Range("AC3").Value = -1615
$ws.Range("AD3").Value = -1.3
$ws.Range("AE3").Value = 1940
$ws.Range("AF3").Value = 1.08
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 32960505

# Row 4
$ws.Range("D4").Value = 2478
$ws.Range("E4").Value = 126
$ws.Range("F4").Value = 126
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = 55
$ws.Range("I4").Value = 52
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 2236
$ws.Range("L4").Value = 1310
$ws.Range("M4").Value = 926
$ws.Range("N4").Value = 877
$ws.Range("O4").Value = 49
$ws.Range("P4").Value = 237
$ws.Range("Q4").Value = 66
$ws.Range("R4").Value = 83
$ws.Range("S4").Value = -196
$ws.Range("T4").Value = 10
$ws.Range("U4").Value = 56
$ws.Range("V4").Value = 1233
$ws.Range("W4").Value = 5.1
$ws.Range("X4").Value = 2.24
$ws.Range("Y4").Value = 6.93
$ws.Range("Z4").Value = 2.36
$ws.Range("AA4").Value = 141.59
$ws.Range("AB4").Value = 175.99
$ws.Range("AC4").Value = 130
$ws.Range("AD4").Value = 15.41
$ws.Range("AE4").Value = 1887
$ws.Range("AF4").Value = 1.06
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 0.5
$ws.Range("AI4").Value = 8.96
$ws.Range("AJ4").Value = 47474590

# Row 5
$ws.Range("D5").Value = 2572
$ws.Range("E5").Value = 128
$ws.Range("F5").Value = 128
$ws.Range("G5").Value = 35
$ws.Range("H5").Value = 20
$ws.Range("I5").Value = 18
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2123
$ws.Range("L5").Value = 1191
$ws.Range("M5").Value = 932
$ws.Range("N5").Value = 881
$ws.Range("O5").Value = 51
$ws.Range("P5").Value = 237
$ws.Range("Q5").Value = 63
$ws.Range("R5").Value = 120
$ws.Range("S5").Value = -151
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = 58
$ws.Range("V5").Value = 1058
$ws.Range("W5").Value = 4.99
$ws.Range("X5").Value = 0.76
$ws.Range("Y5").Value = 2.09
$ws.Range("Z5").Value = 0.9
$ws.Range("AA5").Value = 127.84
$ws.Range("AB5").Value = 179.39
$ws.Range("AC5").Value = 39
$ws.Range("AD5").Value = 29.55
$ws.Range("AE5").Value = 1897
$ws.Range("AF5").Value = 0.6
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 47474590

# Row 6
$ws.Range("D6").Value = 2580
$ws.Range("E6").Value = 86
$ws.Range("F6").Value = 86
$ws.Range("G6").Value = 86
$ws.Range("H6").Value = 70
$ws.Range("I6").Value = 70
$ws.Range("K6").Value = 2201
$ws.Range("L6").Value = 1148
$ws.Range("M6").Value = 1053
$ws.Range("N6").Value = 1003
$ws.Range("P6").Value = 237
$ws.Range("Q6").Value = 32
$ws.Range("R6").Value = -16
$ws.Range("S6").Value = -52
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = 23
$ws.Range("V6").Value = 1005
$ws.Range("W6").Value = 3.33
$ws.Range("X6").Value = 2.73
$ws.Range("Y6").Value = 7.47
$ws.Range("Z6").Value = 3.26
$ws.Range("AA6").Value = 109.05
$ws.Range("AB6").Value = 208.28
$ws.Range("AC6").Value = 148
$ws.Range("AD6").Value = 7.29
$ws.Range("AE6").Value = 2158
$ws.Range("AF6").Value = 0.5
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 1.39
$ws.Range("AI6").Value = 9.9
$ws.Range("AJ6").Value = 47474590

# Clear columns D:AJ for rows 7, 8, 9 (data no longer available for these periods)
$ws.Range("D7:AJ9").ClearContents()
